$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: environment moved to the new pre-production host, new policy number ---
$ws.Range("B2").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("C2").Value = "https://i-preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"

# C2 picks up the same "looks like a hyperlink" look the other URL cells already use (C3:C9)
$ws.Range("C3").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# F2: new policy number, now flagged as quoted text (same border/fill it already had)
$ws.Range("F6").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F2").Value = "04104015957"

# H2: siniestro date updated
$ws.Range("H2").Value = "05/04/2021"

# --- Row 3: reuse the freed-up policy number, new date ---
$ws.Range("F3").Value = "04104015645"
$ws.Range("H3").Value = "22/03/2021"

# --- Sheet view: scroll back to the left edge, move the active selection to I2 ---
$ws.Range("I2").Select()
$excel.ActiveWindow.ScrollColumn = 1
